$d = $word.ActiveDocument

# Title change (appears twice: heading + bold line near the end)
$d.Content.Find.Execute(
    "Play Atlantis Queen Free : Exciting Underwater Slot Game", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Atlantis Queen Free: Game Review & Bonuses", 2
)

# "What we like" bullet list
$d.Content.Find.Execute(
    "Exciting bonus features with potential for big payouts", $true, $false, $false, $false, $false,
    $true, 1, $false, "Pearl Bonus feature offers cash prizes", 2
)

$d.Content.Find.Execute(
    "Attractive graphics and symbols", $true, $false, $false, $false, $false,
    $true, 1, $false, "Atlantis Temple Bonus feature with free spins and multiplier", 2
)

$d.Content.Find.Execute(
    "25 paylines offer excellent chances to win", $true, $false, $false, $false, $false,
    $true, 1, $false, "Underwater world setting and theme", 2
)

# "What we don't like" bullet list
$d.Content.Find.Execute(
    "Base game wins can be small", $true, $false, $false, $false, $false,
    $true, 1, $false, "Limited bonus features", 2
)

$d.Content.Find.Execute(
    "Limited variety in bonus features", $true, $false, $false, $false, $false,
    $true, 1, $false, "No progressive jackpot", 2
)

# Meta description (italic) near the end
$d.Content.Find.Execute(
    "Discover the lost city of Atlantis with the exciting slot game Atlantis Queen. Play for free and experience the thrill of underwater treasures and bonuses.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Read our review of Atlantis Queen and play for free today to enjoy exciting bonuses and cash prizes.", 2
)
